$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "39.504.55"
$ws.Cells.Item(2, 5).Value = "  +1.66%  "

$ws.Cells.Item(3, 4).Value = "2.156.49"
$ws.Cells.Item(3, 5).Value = "  +3.05%  "

$ws.Cells.Item(4, 5).Value = "  +0.06%  "

$ws.Cells.Item(5, 4).Value = "'229.32"
$ws.Cells.Item(5, 5).Value = "  +0.26%  "

$ws.Cells.Item(6, 5).Value = "  +1.04%  "

$ws.Cells.Item(7, 4).Value = "'62.98"
$ws.Cells.Item(7, 5).Value = "  +4.07%  "

$ws.Cells.Item(8, 5).Value = "  +0.01%  "

$ws.Cells.Item(9, 4).Value = "'0.395"
$ws.Cells.Item(9, 5).Value = "  +2.61%  "

$ws.Cells.Item(10, 5).Value = "  +2.84%  "

$ws.Cells.Item(11, 4).Value = "'0.103"
$ws.Cells.Item(11, 5).Value = "  -0.37%  "

$ws.Cells.Item(12, 4).Value = "'16.18"
$ws.Cells.Item(12, 5).Value = "  +7.98%  "

$ws.Cells.Item(13, 4).Value = "2.476.15"
$ws.Cells.Item(13, 5).Value = "  +3.07%  "

$ws.Cells.Item(14, 4).Value = "'22.26"
$ws.Cells.Item(14, 5).Value = "  +1.61%  "

$ws.Cells.Item(15, 4).Value = "'0.822"
$ws.Cells.Item(15, 5).Value = "  +3.35%  "

$ws.Cells.Item(16, 5).Value = "  +1.82%  "

$ws.Cells.Item(17, 4).Value = "2.152.57"
$ws.Cells.Item(17, 5).Value = "  +2.84%  "

$ws.Cells.Item(18, 4).Value = "39.478.94"
$ws.Cells.Item(18, 5).Value = "  +1.92%  "

$ws.Cells.Item(19, 4).Value = "'72.39"
$ws.Cells.Item(19, 5).Value = "  +1.06%  "

$ws.Cells.Item(20, 4).Value = "'6.15"
$ws.Cells.Item(20, 5).Value = "  +1.85%  "

$ws.Cells.Item(21, 5).Value = "  +2.16%  "

$ws.Cells.Item(22, 4).Value = "'228.78"
$ws.Cells.Item(22, 5).Value = "  +0.65%  "

$ws.Cells.Item(23, 5).Value = "  -0.05%  "

$ws.Cells.Item(24, 5).Value = "  +1.08%  "

$ws.Cells.Item(25, 4).Value = "'2.38"
$ws.Cells.Item(25, 5).Value = "  +1.59%  "

$ws.Cells.Item(26, 4).Value = "'9.78"
$ws.Cells.Item(26, 5).Value = "  +2.99%  "

$ws.Cells.Item(27, 4).Value = "'172.69"
$ws.Cells.Item(27, 5).Value = "  +0.81%  "

$ws.Cells.Item(28, 5).Value = "  -0.73%  "

$ws.Cells.Item(29, 5).Value = "  -3.26%  "

$ws.Cells.Item(30, 4).Value = "'19.64"
$ws.Cells.Item(30, 5).Value = "  +2.38%  "

$ws.Cells.Item(31, 4).Value = "'2.57"
$ws.Cells.Item(31, 5).Value = "  +9.19%  "

$ws.Cells.Item(32, 5).Value = "  +1.36%  "

$ws.Cells.Item(33, 4).Value = "'4.63"
$ws.Cells.Item(33, 5).Value = "  +2.88%  "

$ws.Cells.Item(34, 4).Value = "'4.82"
$ws.Cells.Item(34, 5).Value = "  +2.56%  "

$ws.Cells.Item(35, 4).Value = "'7.15"
$ws.Cells.Item(35, 5).Value = "  +10.97%  "

$ws.Cells.Item(36, 4).Value = "'0.0623"
$ws.Cells.Item(36, 5).Value = "  +1.84%  "

$ws.Cells.Item(37, 4).Value = "'2.44"
$ws.Cells.Item(37, 5).Value = "  +2.41%  "

$ws.Cells.Item(38, 4).Value = "'3.57"
$ws.Cells.Item(38, 5).Value = "  -0.33%  "

$ws.Cells.Item(39, 4).Value = "'1.00"
$ws.Cells.Item(39, 5).Value = "  +0.25%  "

$ws.Cells.Item(40, 2).Value = "VeChain"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(40, 4).Value = "'0.0232"
$ws.Cells.Item(40, 5).Value = "  +3.39%  "

$ws.Cells.Item(41, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(41, 4).Value = "'18.16"
$ws.Cells.Item(41, 5).Value = "  -0.02%  "

$ws.Cells.Item(42, 4).Value = "'103.36"

$ws.Cells.Item(43, 4).Value = "1.532.36"
$ws.Cells.Item(43, 5).Value = "  -0.54%  "

$ws.Cells.Item(44, 5).Value = "  +5.63%  "

$ws.Cells.Item(45, 4).Value = "'0.0930"
$ws.Cells.Item(45, 5).Value = "  +0.85%  "

$ws.Cells.Item(46, 2).Value = "FTXToken"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Cells.Item(46, 4).Value = "'4.29"
$ws.Cells.Item(46, 5).Value = "  +4.37%  "

$ws.Cells.Item(47, 2).Value = "ARBITRUM"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(47, 4).Value = "'1.10"
$ws.Cells.Item(47, 5).Value = "  +6.77%  "

$ws.Cells.Item(48, 2).Value = "HuobiToken"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(48, 4).Value = "'2.81"
$ws.Cells.Item(48, 5).Value = "  -0.39%  "

$ws.Cells.Item(49, 2).Value = "FraxShare"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(49, 4).Value = "'7.76"
$ws.Cells.Item(49, 5).Value = "  +1.54%  "

$ws.Cells.Item(50, 4).Value = "2.360.46"
$ws.Cells.Item(50, 5).Value = "  +3.14%  "

$ws.Cells.Item(51, 5).Value = "  +0.03%  "
